$wb = $excel.ActiveWorkbook

# --- Sheet "Hoja1": update the A1 message cell with the new rates ---
$ws1 = $wb.Worksheets.Item("Hoja1")
$newText = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 12.62 = 50511.17 pesos`n✅ 50511.17 pesos = 12.58 = 968.77 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"
$ws1.Range("A1").Value = $newText

# --- Sheet "tasas": update N10, O10, N12, O12 ---
$ws2 = $wb.Worksheets.Item("tasas")
$ws2.Range("N10").Value = 79.25
$ws2.Range("O10").Value = 4003.01
$ws2.Range("N12").Value = 4015
$ws2.Range("O12").Value = 77.005
